$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.043989068301713, 5.13616940762048),
    @(3, 4.363588742201005, 10.35558404259502),
    @(4, 8.533237662149546, 14.95076608011288),
    @(5, 9.466511095318387, 19.9179981284433),
    @(6, 11.64294439302546, 24.66635597833487),
    @(7, 13.30997688524541, 30.28199444573522),
    @(8, 15.51961208849998, 34.97276062348703),
    @(9, 16.58725555717313, 39.6420962801613),
    @(10, 17.41163180709681, 44.16693946149432),
    @(11, 19.91328570273131, 48.89736243076941),
    @(12, 22.0776141222408, 53.63922239201244),
    @(13, 23.47874095452984, 58.90378840985626),
    @(14, 24.12165290875027, 64.26472922535508),
    @(15, 26.38038415588705, 69.10319388144914),
    @(16, 27.0150925882133, 73.78285655139155),
    @(17, 29.05341866659325, 79.19926814031516),
    @(18, 29.88693720240534, 84.15852415259724),
    @(19, 32.16369071755351, 89.79542322471247),
    @(20, 34.85263026938543, 94.62398898104446),
    @(21, 36.39422954453094, 99.02589449258451),
    @(22, 37.28908488591861, 103.841009211044),
    @(23, 40.42136935877834, 108.9693113427763),
    @(24, 44.11142920029371, 114.0948169986551),
    @(25, 46.71339351812423, 118.4872412224238),
    @(26, 47.8062963910805, 123.0897594555255),
    @(27, 51.44206547759828, 127.8493472661304),
    @(28, 55.10058070589977, 132.2891525723353),
    @(29, 56.71930550296261, 137.1185123077073),
    @(30, 57.39409084503254, 141.8763955024542),
    @(31, 58.58551925809187, 147.1494066957908),
    @(32, 60.09049913560703, 151.3209598548814),
    @(33, 60.96758555427176, 155.8254299724525),
    @(34, 63.22136160098436, 161.3283033110014),
    @(35, 64.59861273652717, 166.3857361579106),
    @(36, 65.83225132808828, 171.5993831882853),
    @(37, 67.49252721370647, 176.619382611278),
    @(38, 70.25875877705113, 181.2819182761598),
    @(39, 72.35330502517238, 185.6249571936228),
    @(40, 75.28668362260612, 191.1095766902215),
    @(41, 77.7253803218552, 195.9149614591009),
    @(42, 79.15234668990063, 200.9886025623993),
    @(43, 85.02180152514413, 206.0949967755045),
    @(44, 86.25750359926629, 210.7021400700368),
    @(45, 89.39865241737415, 214.8790465407601),
    @(46, 91.67356707135421, 219.9688075817991),
    @(47, 92.57749897978481, 225.3474909829673),
    @(48, 94.72152651894899, 230.1766474558695),
    @(49, 96.96846696674989, 235.305886518541),
    @(50, 99.1356606236715, 240.0804058657262),
)

foreach ($row in $data) {
    $r = $row[0]
    $bVal = $row[1]
    $cVal = $row[2]
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
}
